$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 514, shifting existing rows 514:535 down to 515:536
$ws.Rows.Item(514).Insert()

# Populate the newly inserted row 514 with the new data
$ws.Cells.Item(514, 1).Value = 5
$ws.Cells.Item(514, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(514, 3).Value = "Maule"
$ws.Cells.Item(514, 4).Value = 45267
$ws.Cells.Item(514, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(514, 5).Value = 7
$ws.Cells.Item(514, 6).Value = 100112045
$ws.Cells.Item(514, 7).Value = "Zapallo"
$ws.Cells.Item(514, 8).Value = "Paine"
$ws.Cells.Item(514, 9).Value = "1a (guarda)"
$ws.Cells.Item(514, 10).Value = 1500
$ws.Cells.Item(514, 11).Value = 900
$ws.Cells.Item(514, 12).Value = 900
$ws.Cells.Item(514, 13).Value = 900
$ws.Cells.Item(514, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(514, 15).Value = "Región del Maule"
$ws.Cells.Item(514, 16).Value = 900
$ws.Cells.Item(514, 17).Value = 1
$ws.Cells.Item(514, 18).Value = "Hortaliza"
